$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray formatted-but-empty cell left over at J15 from earlier edits
$ws.Range("J15").Clear()

# New "Estados / LOAD / STOP / DIRECCION" table data (rows 31-36)
$ws.Range("A31").Value = "Estados"
$ws.Range("B31").Value = "LOAD"
$ws.Range("C31").Value = "STOP"
$ws.Range("D31").Value = "DIRECCIÓN"
$ws.Range("E31").Value = "Columna5"
$ws.Range("F31").Value = "Columna6"
$ws.Range("G31").Value = "Columna7"
$ws.Range("H31").Value = "Columna8"
$ws.Range("I31").Value = "Columna9"

$ws.Range("A32").Value = "Estado 1"
$ws.Range("B32").Value = 1
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = "X"

$ws.Range("A33").Value = "Estado 2"
$ws.Range("B33").Value = 0
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = "X"

$ws.Range("A34").Value = "Estado 3"
$ws.Range("B34").Value = 0
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = "X"

$ws.Range("A35").Value = "Estado 4"
$ws.Range("B35").Value = 0
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = "X"

$ws.Range("A36").Value = "Estado 5"
$ws.Range("B36").Value = 0
$ws.Range("C36").Value = 1
$ws.Range("D36").Value = "X"

# Turn the new range into an Excel table ("Tabla7"), same style as the others
$lo = $ws.ListObjects.Add(1, $ws.Range("A31:I36"), 0, 1)
$lo.Name = "Tabla7"
$lo.TableStyle = "TableStyleMedium2"

# Match the new selection / active cell the author ended up on
$ws.Range("D31").Select() | Out-Null
